# New GNG latency code (works for both L/R)
#
# The hit-rate / false-alarm-rate calculations now use a fixed trial
# denominator of 30 (instead of the previous per-row trial counts), and
# the lever-press latency columns (Go / NoGo) are recomputed by the new
# latency code using the raw press timestamps.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$trialDenom = 30

# New lever-press latency values (Go = column I, NoGo = column J),
# produced by the new GNG latency code, keyed by worksheet row.
$latency = @{
    2  = @(1.945, 2.257)
    3  = @(2.219, 1.798)
    4  = @(2.198, 2.775)
    5  = @(2.893, 2.099)
    6  = @(1.553, 2.352)
    7  = @(1.746, 2.117)
    8  = @(2.679, 2.967)
    9  = @(1.359, 1.623)
    10 = @(2.087, 2.111)
    11 = @(2.186, 2.251)
    12 = @(1.963, 1.806)
    13 = @(2.576, 2.841)
    14 = @(1.53,  1.85)
    15 = @(2.493, 2.268)
    16 = @(2.119, 2.899)
    17 = @(2.44,  2.802)
    18 = @(1.206, 1.344)
    19 = @(3.029, 3.461)
    20 = @(3.072, 1.624)
    21 = @(1.353, 1.62)
    22 = @(2.891, 3.04)
    23 = @(1.696, 2.078)
    24 = @(2.211, 2.28)
    25 = @(2.176, 2.37)
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($row = 2; $row -le $lastRow; $row++) {
    $goSuccess = $ws.Cells.Item($row, 5).Value()    # column E - Successful Go Trials
    $noGoSuccess = $ws.Cells.Item($row, 6).Value()  # column F - Successful NoGo Trials

    # Hit Rate (column G) and False Alarm Rate (column H), now computed
    # against the fixed 30-trial denominator.
    $ws.Cells.Item($row, 7).Value = ($goSuccess / $trialDenom) * 100
    $ws.Cells.Item($row, 8).Value = (($trialDenom - $noGoSuccess) / $trialDenom) * 100

    if ($latency.ContainsKey($row)) {
        $vals = $latency[$row]
        $ws.Cells.Item($row, 9).Value = $vals[0]   # column I - Lever Press Latency Go
        $ws.Cells.Item($row, 10).Value = $vals[1]  # column J - Lever Press Latency NoGo
    }
}
